$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as TEXT, preserving strings that look like numbers
# (e.g. "584.66", "1.00") instead of letting them get auto-coerced into
# numeric cells. We force the Text number format for the write, then clear
# the formatting override again so the cell's style matches the original
# (unstyled) cell.
function Set-TextCell($addr, $value) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.ClearFormats()
}

# Row 2 - Bitcoin
Set-TextCell "D2" "62.800.30"
Set-TextCell "E2" "  +4.86%  "

# Row 3 - Ethereum
Set-TextCell "D3" "3.111.76"
Set-TextCell "E3" "  +2.98%  "

# Row 5 - BNB
Set-TextCell "D5" "584.66"
Set-TextCell "E5" "  +3.25%  "

# Row 6 - Solana
Set-TextCell "D6" "144.60"
Set-TextCell "E6" "  +2.39%  "

# Row 8 - LidoStakedEther
Set-TextCell "D8" "3.104.78"
Set-TextCell "E8" "  +3.11%  "

# Row 9 - XRP
Set-TextCell "D9" "0.530"
Set-TextCell "E9" "  +1.52%  "

# Row 10 - Dogecoin
Set-TextCell "E10" "  +10.55%  "

# Row 11 - Toncoin
Set-TextCell "D11" "5.74"
Set-TextCell "E11" "  +7.98%  "

# Row 12 - Cardano
Set-TextCell "E12" "  +1.32%  "

# Row 13 - ShibaInu
Set-TextCell "E13" "  +6.02%  "

# Row 14 - Avalanche
Set-TextCell "D14" "35.51"
Set-TextCell "E14" "  +3.68%  "

# Row 15 - TRON
Set-TextCell "E15" "  -0.13%  "

# Row 16 - WrappedliquidstakedEther2.0
Set-TextCell "D16" "3.626.07"
Set-TextCell "E16" "  +2.99%  "

# Row 17 - Polkadot
Set-TextCell "E17" "  -0.69%  "

# Row 18 - WrappedEther
Set-TextCell "D18" "3.108.31"
Set-TextCell "E18" "  +2.92%  "

# Row 19 - WrappedBTC
Set-TextCell "D19" "62.738.91"
Set-TextCell "E19" "  +4.81%  "

# Row 20 - BitcoinCash
Set-TextCell "D20" "465.07"
Set-TextCell "E20" "  +5.79%  "

# Row 21 - Chainlink
Set-TextCell "E21" "  +2.64%  "

# Row 22 - Polygon
Set-TextCell "E22" "  +0.74%  "

# Row 23 - Uniswap
Set-TextCell "D23" "7.55"
Set-TextCell "E23" "  +5.76%  "

# Row 24 - InternetComputer(DFINITY)
Set-TextCell "D24" "13.33"
Set-TextCell "E24" "  -0.36%  "

# Row 25 - Litecoin
Set-TextCell "D25" "82.12"
Set-TextCell "E25" "  +1.56%  "

# Row 26 - Dai
Set-TextCell "D26" "1.00"
Set-TextCell "E26" "  -0.03%  "

# Row 27 - ImmutableX
Set-TextCell "E27" "  -0.04%  "

# Row 28 - PancakeSwap
Set-TextCell "E28" "  +4.87%  "

# Row 29 - FirstDigitalUSD
Set-TextCell "E29" "  -0.05%  "

# Row 30 - RenderToken
Set-TextCell "E30" "  +5.17%  "

# Row 31 - NEARProtocol
Set-TextCell "D31" "6.82"
Set-TextCell "E31" "  +8.29%  "

# Row 32 & 33 - EthereumClassic / Hedera swap places
Set-TextCell "B32" "Hedera"
Set-TextCell "C32" "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextCell "D32" "0.111"
Set-TextCell "E32" "  +7.16%  "

Set-TextCell "B33" "EthereumClassic"
Set-TextCell "C33" "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextCell "D33" "26.91"
Set-TextCell "E33" "  +3.12%  "

# Row 34 - PEPE
Set-TextCell "D34" "0.0₃0852"
Set-TextCell "E34" "  +6.95%  "

# Row 35 - Stacks
Set-TextCell "E35" "  +11.27%  "

# Row 36 - Mantle
Set-TextCell "E36" "  +3.50%  "

# Row 37 - Filecoin
Set-TextCell "E37" "  +1.86%  "

# Row 38 - dogwifhat
Set-TextCell "D38" "3.26"
Set-TextCell "E38" "  +15.62%  "

# Row 39 - OKB
Set-TextCell "D39" "50.92"
Set-TextCell "E39" "  +3.45%  "

# Row 40 - Bittensor
Set-TextCell "D40" "433.05"
Set-TextCell "E40" "  +6.50%  "

# Row 41 - Cosmos
Set-TextCell "D41" "8.76"
Set-TextCell "E41" "  +1.06%  "

# Row 42 - Maker
Set-TextCell "D42" "2.937.03"
Set-TextCell "E42" "  +5.53%  "

# Row 43 - VeChain
Set-TextCell "E43" "  +4.01%  "

# Row 44 - TheGraph
Set-TextCell "E44" "  +8.79%  "

# Row 45 - Kaspa
Set-TextCell "E45" "  +3.03%  "

# Row 46 - Fetch.AI
Set-TextCell "E46" "  +6.46%  "

# Row 47 - Arweave
Set-TextCell "D47" "35.36"
Set-TextCell "E47" "  +4.05%  "

# Row 49 - Monero
Set-TextCell "D49" "123.36"
Set-TextCell "E49" "  +0.39%  "

# Row 50 - Stellar
Set-TextCell "E50" "  +0.33%  "

# Row 51 - InjectiveProtocol
Set-TextCell "D51" "24.75"
Set-TextCell "E51" "  +4.39%  "
